$wb = $excel.ActiveWorkbook

# --- Sheet "Games": append the newly-played game as row 40 ---
$games = $wb.Worksheets.Item("Games")

$games.Range("A40").Value = 39
$games.Range("B40").Value = 45305
$games.Range("B40").NumberFormat = "YYYY-MM-DD"
$games.Range("C40").Value = 2
$games.Range("D40").Value = 127
$games.Range("E40").Value = 94.3
$games.Range("F40").Value = 0.6909999999999999
$games.Range("G40").Value = 12
$games.Range("H40").Value = 22.6
$games.Range("I40").Value = 0.289
$games.Range("J40").Value = 134.6
$games.Range("K40").Value = "POR"
$games.Range("L40").Value = 116
$games.Range("M40").Value = 0.5620000000000001
$games.Range("N40").Value = 11.6
$games.Range("O40").Value = 26.1
$games.Range("P40").Value = 0.18
$games.Range("Q40").Value = 123
$games.Range("R40").Value = 0
$games.Range("S40").Value = 1

# --- Sheet "Next": the game above is no longer upcoming, so drop its row ---
# (row 2 = 2023-12-16 vs POR) and let every later fixture shift up by one.
$next = $wb.Worksheets.Item("Next")
$next.Rows.Item(2).Delete()
